$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so values like
# "1.002" or "0.000007860" are not silently coerced into numbers, then
# reset the style back to Normal so no stray per-cell style index remains.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '30.702.67'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '1.889.04'
$ws.Range('E3').Value = '  -1.70%  '
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').Value = '235.55'
$ws.Range('E5').Value = '  -4.46%  '
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.22%  '
$ws.Range('D7').Value = '0.4885'
$ws.Range('E7').Value = '  -2.64%  '
$ws.Range('D8').Value = '0.2899'
$ws.Range('E8').Value = '  -4.10%  '
$ws.Range('D9').Value = '0.06669'
$ws.Range('E9').Value = '  -4.39%  '
$ws.Range('D10').Value = '1.888.53'
$ws.Range('E10').Value = '  -1.68%  '
$ws.Range('D11').Value = '16.65'
$ws.Range('E11').Value = '  -3.11%  '
$ws.Range('D12').Value = '0.07239'
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('D13').Value = '89.07'
$ws.Range('E13').Value = '  -3.96%  '
$ws.Range('D14').Value = '5.005'
$ws.Range('E14').Value = '  -2.30%  '
$ws.Range('D15').Value = '0.6648'
$ws.Range('E15').Value = '  -3.19%  '
$ws.Range('D16').Value = '30.654.22'
$ws.Range('E16').Value = '  -1.56%  '
$ws.Range('D17').Value = '0.000007860'
$ws.Range('E17').Value = '  -3.16%  '
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').Value = '12.98'
$ws.Range('E19').Value = '  -3.64%  '
$ws.Range('D20').Value = '2.126.22'
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('E21').Value = '  +0.38%  '
$ws.Range('D22').Value = '4.735'
$ws.Range('E22').Value = '  -3.47%  '
$ws.Range('D23').Value = '188.54'
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '6.057'
$ws.Range('E24').Value = '  -1.01%  '
$ws.Range('D25').Value = '9.304'
$ws.Range('E25').Value = '  -1.27%  '
$ws.Range('D26').Value = '158.09'
$ws.Range('E26').Value = '  +3.59%  '
$ws.Range('D27').Value = '18.27'
$ws.Range('E27').Value = '  -1.16%  '
$ws.Range('D28').Value = '1.832'
$ws.Range('E28').Value = '  -6.85%  '
$ws.Range('D29').Value = '1.402'
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').Value = '4.264'
$ws.Range('E30').Value = '  -2.67%  '
$ws.Range('D31').Value = '0.09021'
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').Value = '3.942'
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('D33').Value = '0.05189'
$ws.Range('E33').Value = '  -1.54%  '
$ws.Range('D34').Value = '0.7318'
$ws.Range('E34').Value = '  -3.21%  '
$ws.Range('D35').Value = '1.083'
$ws.Range('E35').Value = '  -5.98%  '
$ws.Range('D36').Value = '2.693'
$ws.Range('D37').Value = '0.01825'
$ws.Range('E37').Value = '  -6.24%  '
$ws.Range('D38').Value = '2.669'
$ws.Range('E38').Value = '  -3.00%  '
$ws.Range('D39').Value = '0.9217'
$ws.Range('E39').Value = '  -2.38%  '
$ws.Range('D40').Value = '2.047'
$ws.Range('E40').Value = '  -7.58%  '
$ws.Range('D41').Value = '0.4413'
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').Value = '104.78'
$ws.Range('E42').Value = '  -1.20%  '
$ws.Range('D43').Value = '0.9999'
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('D44').Value = '5.731'
$ws.Range('E44').Value = '  -3.92%  '
$ws.Range('D45').Value = '0.1345'
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '7.294'
$ws.Range('E46').Value = '  -7.89%  '
$ws.Range('D47').Value = '0.4093'
$ws.Range('E47').Value = '  +4.60%  '
$ws.Range('D48').Value = '0.05830'
$ws.Range('E48').Value = '  -0.96%  '
$ws.Range('D49').Value = '8.679'
$ws.Range('E49').Value = '  +0.06%  '
$ws.Range('D50').Value = '1.409'
$ws.Range('E50').Value = '  +0.13%  '
$ws.Range('D51').Value = '33.23'
$ws.Range('E51').Value = '  -0.73%  '

$ws.Range("D2:E51").Style = "Normal"
